$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("A 2025-2027")
$ws2 = $wb.Worksheets.Item("A 2025-2026")

# --- Fill in the next three upcoming matches on the "A 2025-2026" sheet ---
# (our club "Lempo A2" in column C, opponent in column D)
$ws2.Activate()

# Row 6 first, then row 7, then row 5 so the new shared-string entries are
# created in the same order the source workbook uses (VaLePa, OsVa PU17, Jou).
$ws2.Range("C6").Value = "Lempo A2"
$ws2.Range("D6").Value = "VaLePa"

$ws2.Range("C7").Value = "Lempo A2"
$ws2.Range("D7").Value = "OsVa PU17"

$ws2.Range("C5").Value = "Lempo A2"
$ws2.Range("D5").Value = "Jou"

# --- Update the saved selection / active sheet state ---
# Sheet "A 2025-2027" keeps its own remembered selection (no longer the active tab)
$ws1.Activate()
$ws1.Range("H12").Select()

# "A 2025-2026" ends up as the active (visible) sheet with its own selection
$ws2.Activate()
$ws2.Range("F14").Select()
